$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove existing hyperlinks (and their relationships) so we can rebuild
# them cleanly against the refreshed row data below.
$ws.Hyperlinks.Delete()

# Row 2
$ws.Cells.Item(2, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(2, 2).Value = "大手SIer等のAIソリューション開発・導入を支援してくださるエンジニア・PM募集"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5450158"
$ws.Cells.Item(2, 7).Value = 368
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai ◆開発"
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://www.lancers.jp/work/detail/5450158", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5450158") | Out-Null

# Row 3
$ws.Cells.Item(3, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(3, 2).Value = "EC×AIプロダクト/業務改善リード"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5450024"
$ws.Cells.Item(3, 7).Value = 338
$ws.Cells.Item(3, 8).Value = "🔥AI,Ai ◇業務改善"
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://www.lancers.jp/work/detail/5450024", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5450024") | Out-Null

# Row 4
$ws.Cells.Item(4, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(4, 2).Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Cells.Item(4, 7).Value = 243
$ws.Cells.Item(4, 8).Value = "🔥API ◆ツール"
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://www.lancers.jp/work/detail/5217096", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5217096") | Out-Null

# Row 5
$ws.Cells.Item(5, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(5, 2).Value = "【自動化】Webサービス更新ツール開発(200アカウント管理)"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5448409"
$ws.Cells.Item(5, 7).Value = 230
$ws.Cells.Item(5, 8).Value = "◆ツール,開発 ◇管理"
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://www.lancers.jp/work/detail/5448409", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5448409") | Out-Null

# Row 6
$ws.Cells.Item(6, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(6, 2).Value = "【急募】pythonのコードのMac環境用インストーラー作成(Windows版は作成済)"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5442448"
$ws.Cells.Item(6, 7).Value = 190
$ws.Cells.Item(6, 8).Value = "🔥Python"
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://www.lancers.jp/work/detail/5442448", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5442448") | Out-Null

# Row 7
$ws.Cells.Item(7, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(7, 2).Value = "【バイナリ解析 / 逆コンパイル】EPCデータ解析ツール開発"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5450504"
$ws.Cells.Item(7, 7).Value = 128
$ws.Cells.Item(7, 8).Value = "◆ツール,開発"
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), "https://www.lancers.jp/work/detail/5450504", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5450504") | Out-Null

# Row 8
$ws.Cells.Item(8, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(8, 2).Value = "【バイナリ解析 / 逆コンパイル】EPCデータ解析ツール開発(継続依頼あり・高単価)"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5449973"
$ws.Cells.Item(8, 7).Value = 128
$ws.Cells.Item(8, 8).Value = "◆ツール,開発"
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), "https://www.lancers.jp/work/detail/5449973", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5449973") | Out-Null

# Row 9
$ws.Cells.Item(9, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(9, 2).Value = "報酬計算・源泉所得税・支払明細自動化システムの構築(Excel/スプレッドシートベース)"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5450283"
$ws.Cells.Item(9, 7).Value = 98
$ws.Cells.Item(9, 8).Value = "◆自動化"
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), "https://www.lancers.jp/work/detail/5450283", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5450283") | Out-Null

# Row 10
$ws.Cells.Item(10, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(10, 2).Value = "JavaScriptをスクラッチで対応できるパートナー募集"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5450393"
$ws.Cells.Item(10, 7).Value = 78
$ws.Cells.Item(10, 8).Value = "★Java"
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), "https://www.lancers.jp/work/detail/5450393", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5450393") | Out-Null

# Row 11
$ws.Cells.Item(11, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(11, 2).Value = "【単発/Stripeサブスクリプション実装】Laravel プラットフォーム開発エンジニア募集"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5449939"
$ws.Cells.Item(11, 7).Value = 75
$ws.Cells.Item(11, 8).Value = "◆開発"
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), "https://www.lancers.jp/work/detail/5449939", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5449939") | Out-Null

# Row 12
$ws.Cells.Item(12, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(12, 2).Value = "Excel VBA一括自動処理ツール作成(データ転記・分類・置換・NGチェック)【エクセルマクロ】"
$ws.Cells.Item(12, 3).Value = "システム開発"
$ws.Cells.Item(12, 4).Value = "1,000 ~ 5,000 円 / 固定"
$ws.Cells.Item(12, 5).Value = "期限情報なし"
$ws.Cells.Item(12, 6).Value = "https://www.lancers.jp/work/detail/5450139"
$ws.Cells.Item(12, 7).Value = 65
$ws.Cells.Item(12, 8).Value = "◆ツール"
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), "https://www.lancers.jp/work/detail/5450139", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5450139") | Out-Null

# Row 13
$ws.Cells.Item(13, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(13, 2).Value = "初回 【案件】Win2008(PHP5.3)→ Linux(LAMP)へのレガシー調査と移行"
$ws.Cells.Item(13, 3).Value = "システム開発"
$ws.Cells.Item(13, 4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(13, 5).Value = "期限情報なし"
$ws.Cells.Item(13, 6).Value = "https://www.lancers.jp/work/detail/5449999"
$ws.Cells.Item(13, 7).Value = 40
$ws.Cells.Item(13, 8).Value = "○PHP"
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), "https://www.lancers.jp/work/detail/5449999", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5449999") | Out-Null

# Row 14
$ws.Cells.Item(14, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(14, 2).Value = "【急募】古いPHPとPerlプログラムのアップデート依頼"
$ws.Cells.Item(14, 3).Value = "システム開発"
$ws.Cells.Item(14, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(14, 5).Value = "期限情報なし"
$ws.Cells.Item(14, 6).Value = "https://www.lancers.jp/work/detail/5440861"
$ws.Cells.Item(14, 7).Value = 33
$ws.Cells.Item(14, 8).Value = "○PHP"
$ws.Hyperlinks.Add($ws.Cells.Item(14, 6), "https://www.lancers.jp/work/detail/5440861", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5440861") | Out-Null

# Row 15
$ws.Cells.Item(15, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(15, 2).Value = "【急募】企業のセキュリティ対策を担うエンジニア募集"
$ws.Cells.Item(15, 3).Value = "システム開発"
$ws.Cells.Item(15, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(15, 5).Value = "期限情報なし"
$ws.Cells.Item(15, 6).Value = "https://www.lancers.jp/work/detail/5450345"
$ws.Cells.Item(15, 7).Value = 25
$ws.Hyperlinks.Add($ws.Cells.Item(15, 6), "https://www.lancers.jp/work/detail/5450345", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5450345") | Out-Null

# Row 16
$ws.Cells.Item(16, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(16, 2).Value = "【急募】社内システム保守運用・社内スタッフ教育まで依頼できる方を探しています!"
$ws.Cells.Item(16, 3).Value = "システム開発"
$ws.Cells.Item(16, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(16, 5).Value = "期限情報なし"
$ws.Cells.Item(16, 6).Value = "https://www.lancers.jp/work/detail/5449609"
$ws.Cells.Item(16, 7).Value = 25
$ws.Hyperlinks.Add($ws.Cells.Item(16, 6), "https://www.lancers.jp/work/detail/5449609", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5449609") | Out-Null

# Row 17
$ws.Cells.Item(17, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(17, 2).Value = "注目 限定公開 PR 限定公開の仕事"
$ws.Cells.Item(17, 3).Value = "システム開発"
$ws.Cells.Item(17, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(17, 5).Value = "期限情報なし"
$ws.Cells.Item(17, 6).Value = "https://www.lancers.jp/work/detail/5450323"
$ws.Cells.Item(17, 7).Value = 13
$ws.Hyperlinks.Add($ws.Cells.Item(17, 6), "https://www.lancers.jp/work/detail/5450323", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5450323") | Out-Null

# Row 18
$ws.Cells.Item(18, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(18, 2).Value = "自社カレンダーとGoogleカレンダーの連携エキスパート募集"
$ws.Cells.Item(18, 3).Value = "システム開発"
$ws.Cells.Item(18, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(18, 5).Value = "期限情報なし"
$ws.Cells.Item(18, 6).Value = "https://www.lancers.jp/work/detail/5450296"
$ws.Cells.Item(18, 7).Value = 13
$ws.Hyperlinks.Add($ws.Cells.Item(18, 6), "https://www.lancers.jp/work/detail/5450296", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5450296") | Out-Null

# Row 19
$ws.Cells.Item(19, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(19, 2).Value = "X(旧ツイッター)自動ログインについて"
$ws.Cells.Item(19, 3).Value = "システム開発"
$ws.Cells.Item(19, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(19, 5).Value = "期限情報なし"
$ws.Cells.Item(19, 6).Value = "https://www.lancers.jp/work/detail/5449817"
$ws.Cells.Item(19, 7).Value = 13
$ws.Hyperlinks.Add($ws.Cells.Item(19, 6), "https://www.lancers.jp/work/detail/5449817", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5449817") | Out-Null

# Row 20
$ws.Cells.Item(20, 1).Value = "2025-12-09 18:23:53"
$ws.Cells.Item(20, 2).Value = "【アカウント復活】削除したxのアカウントを再生させたい!"
$ws.Cells.Item(20, 3).Value = "システム開発"
$ws.Cells.Item(20, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(20, 5).Value = "期限情報なし"
$ws.Cells.Item(20, 6).Value = "https://www.lancers.jp/work/detail/5449948"
$ws.Cells.Item(20, 7).Value = 10
$ws.Hyperlinks.Add($ws.Cells.Item(20, 6), "https://www.lancers.jp/work/detail/5449948", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5449948") | Out-Null

